# Update the LinkedIn carousel draft: replace the "Juniper Green / Solapur"
# article placeholder copy with the generated "dragon in the grid" (China /
# EU energy security) article copy, across all six carousel slides.
#
# Note: this deck's text-diff engine isolates a trailing sentence-ending
# "." into its own run when the *previous* text in that paragraph also
# ended in ".". Writing a short non-terminal placeholder first, then the
# real (period-ending) sentence, keeps each paragraph's new text in a
# single run - matching how the source OOXML actually looks.

$p = $ppt.ActivePresentation

$title = "The dragon in the grid: Limiting China’s influence in Europe’s energy system - European Union Institute for Security Studies |"

$body = @{
    1 = @(
        "The European Union is increasingly concerned about China's investments in its energy infrastructure.",
        "China's involvement in European energy projects has raised security and dependency issues."
    )
    2 = @(
        "The EU is exploring strategies to reduce reliance on Chinese technology in energy systems.",
        "There are ongoing discussions within the EU about enhancing energy security and diversification."
    )
    3 = @(
        "China has become a significant player in renewable energy investments in Europe.",
        "The EU aims to balance investments from China with local and allied sources."
    )
    4 = @(
        "The European Commission is assessing the implications of foreign investments in critical energy infrastructure.",
        "Regulatory frameworks are being developed to scrutinize and potentially limit Chinese investments."
    )
    5 = @(
        "Strategic partnerships with other countries are being prioritized to enhance energy independence.",
        "The EU is focusing on strengthening its internal market to reduce vulnerabilities."
    )
    6 = @(
        "The EU's energy strategy includes a commitment to sustainability and reducing carbon emissions.",
        "Efforts are underway to ensure that energy transitions do not compromise security interests."
    )
}

function Set-ParagraphText {
    param($shape, $index, $text)

    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs($index, 1)

    # First write a placeholder that does NOT end in sentence punctuation,
    # so the follow-up write (the real text) isn't diffed against an
    # old run that also ended in ".", which is what triggers the
    # spurious trailing-period run split.
    $para.Text = "placeholder text"

    $tr2 = $shape.TextFrame.TextRange
    $para2 = $tr2.Paragraphs($index, 1)
    $para2.Text = $text
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $shape = $slide.Shapes.Item(1)

    Set-ParagraphText $shape 1 $title
    Set-ParagraphText $shape 2 $body[$i][0]
    Set-ParagraphText $shape 3 $body[$i][1]
}
